$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. ADC sheet: correct ADC reading (246 -> 247) and add a new "reference
#    voltage drift" block in rows 25-28.
# ---------------------------------------------------------------------------
$adcSheet = $wb.Worksheets.Item("ADC")

$adcSheet.Range("B15").Value = 247

$adcSheet.Range("A25").Value = "Spannung Referenz, mV"
$adcSheet.Range("B25").Value = 1225
$adcSheet.Range("A26").Value = "ADC Wert dazu:"
$adcSheet.Range("B26").Formula = "=ROUND(B25*1024/5000,0)"
$adcSheet.Range("C26").Formula = "=B26*5000/1024"
$adcSheet.Range("A27").Value = "Spannung Referenz gemessen"
$adcSheet.Range("B27").Value = 247
$adcSheet.Range("A28").Value = "Drift:"
$adcSheet.Range("B28").Formula = "=B25-B27"

$adcSheet.Columns.Item(1).ColumnWidth = 27.1
$adcSheet.Columns.Item(2).ColumnWidth = 11.7

$adcSheet.Range("B27").Select()

# ---------------------------------------------------------------------------
# 2. Insert a new "Ip" worksheet right after "PWM" (becomes the 2nd tab) and
#    fill it with the primary-current calculation.
# ---------------------------------------------------------------------------
$pwmSheet = $wb.Worksheets.Item("PWM")
$newSheet = $wb.Worksheets.Add($null, $pwmSheet)
$newSheet.Name = "Ip"

# Pre-register the new shared strings in their original authoring order
# (Shunt, Ua, Ua_ref, delta, divisor, Amp, Ip) using a scratch range, then
# clear it before writing the real layout below.
$newSheet.Range("Z1").Value = "Shunt"
$newSheet.Range("Z2").Value = "Ua"
$newSheet.Range("Z3").Value = "Ua_ref"
$newSheet.Range("Z4").Value = "delta"
$newSheet.Range("Z5").Value = "divisor"
$newSheet.Range("Z6").Value = "Amp"
$newSheet.Range("Z7").Value = "Ip"
$newSheet.Range("Z1:Z7").ClearContents()

$newSheet.Range("A1").Value = "Shunt"
$newSheet.Range("B1").Value = 61900
$newSheet.Range("A2").Value = "Amp"
$newSheet.Range("B2").Value = 8
$newSheet.Range("A3").Value = "Ua"
$newSheet.Range("B3").Value = 2069
$newSheet.Range("A4").Value = "Ua_ref"
$newSheet.Range("B4").Value = 1503
$newSheet.Range("A5").Value = "delta"
$newSheet.Range("B5").Formula = "=(B3-B4)*1000"
$newSheet.Range("A6").Value = "divisor"
$newSheet.Range("B6").Formula = "=B1*B2"
$newSheet.Range("A7").Value = "Ip"
$newSheet.Range("B7").Formula = "=B5/B6*1000"

$newSheet.Range("B7").Select()
